$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''307.60'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''0.76%'
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''38.58'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''8.41%'
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = '''5.106'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '''1.14%'
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''0.08115'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''1.30%'
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''1.972'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''5.18%'
$ws.Range('E6').Style = 'Normal'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').Value = '''4.184'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''1.26%'
$ws.Range('E7').Style = 'Normal'
$ws.Range('B8').Value = 'KuCoinToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D8').Value = '''7.945'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''2.04%'
$ws.Range('E8').Style = 'Normal'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = '''0.9291'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''0.63%'
$ws.Range('E9').Style = 'Normal'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = '''0.1431'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''11.08%'
$ws.Range('E10').Style = 'Normal'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = '''0.1956'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''3.04%'
$ws.Range('E11').Style = 'Normal'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '''0.09061'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''-0.79%'
$ws.Range('E12').Style = 'Normal'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '''0.03504'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''2.57%'
$ws.Range('E13').Style = 'Normal'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '''0.09821'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''-0.38%'
$ws.Range('E14').Style = 'Normal'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '''0.001405'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''-0.03%'
$ws.Range('E15').Style = 'Normal'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '''0.006132'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''-1.30%'
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '''3.733'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''-2.92%'
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''3.446'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''2.80%'
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = '''1.30%'
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''0.1294'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''-4.00%'
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''4.792'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''-0.54%'
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''0.2455'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''6.35%'
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''0.04354'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''-1.65%'
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = '''-1.02%'
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = '''-1.04%'
$ws.Range('E25').Style = 'Normal'
$ws.Range('D27').Value = '''0.0001301'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''3.94%'
$ws.Range('E27').Style = 'Normal'
$ws.Range('D39').Value = '''0.02087'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''7.74%'
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''0.05108'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''-1.69%'
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''0.007475'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''-1.04%'
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = '''-0.07%'
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''0.1358'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''0.62%'
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''0.002131'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''-1.90%'
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''0.009271'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''-3.74%'
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''0.00006258'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''1.26%'
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = '''0.04%'
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''0.003027'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E49').Value = '''-3.56%'
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = '''0.04%'
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = '''0.04%'
$ws.Range('E51').Style = 'Normal'
